$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New handoff/handback report run: drop the "ea39ce26-...md" row (row 3) from
# every sheet, and stamp the remaining (91bbbb67-...md) row with the new
# handoff/handback timestamps from the fresh run.
# ---------------------------------------------------------------------------

$zhHandoffDt  = "2016-03-24 10:23:10"
$zhHandbackDt = "2016-03-24 10:23:53"
$deHandoffDt  = "2016-03-24 10:23:19"
$deHandbackDt = "2016-03-24 10:24:11"

# --- Overview sheet: just file name + zh-cn / de-de status columns --------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Rows.Item(3).Delete()
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/14fcd17033617a21dc8936b2ebe8a2e90f707523/e2e/91bbbb67-038f-4402-9f8c-f2af61374c0d.md", "", "", "91bbbb67-038f-4402-9f8c-f2af61374c0d.md") | Out-Null

# --- zh-cn sheet ------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Rows.Item(3).Delete()
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/14fcd17033617a21dc8936b2ebe8a2e90f707523/e2e/91bbbb67-038f-4402-9f8c-f2af61374c0d.md", "", "", "91bbbb67-038f-4402-9f8c-f2af61374c0d.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cc9f2f0379660b67968ac455bcd3d86e9aaa1106/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/91bbbb67-038f-4402-9f8c-f2af61374c0d.c31ca651c21b93c2e8ff5d2f0b01ba6f1462c3e1.zh-cn.xlf", "", "", "91bbbb67-038f-4402-9f8c-f2af61374c0d.c31ca651c21b93c2e8ff5d2f0b01ba6f1462c3e1.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/ebdd39d73b7e9cee9bc8f80a36c5c44417cd160a/e2e/91bbbb67-038f-4402-9f8c-f2af61374c0d.md", "", "", "91bbbb67-038f-4402-9f8c-f2af61374c0d.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3cc25921b6c7118432ffe64c67e0a21816cc2b46/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/91bbbb67-038f-4402-9f8c-f2af61374c0d.c31ca651c21b93c2e8ff5d2f0b01ba6f1462c3e1.zh-cn.xlf", "", "", "91bbbb67-038f-4402-9f8c-f2af61374c0d.c31ca651c21b93c2e8ff5d2f0b01ba6f1462c3e1.zh-cn.xlf") | Out-Null

$wsZh.Range("E2").Value = $zhHandoffDt
$wsZh.Range("H2").Value = $zhHandbackDt

# --- de-de sheet ------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Rows.Item(3).Delete()
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/14fcd17033617a21dc8936b2ebe8a2e90f707523/e2e/91bbbb67-038f-4402-9f8c-f2af61374c0d.md", "", "", "91bbbb67-038f-4402-9f8c-f2af61374c0d.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a59f9e825a8e627d2a5f625bddd3fcd64fc3913/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/91bbbb67-038f-4402-9f8c-f2af61374c0d.c31ca651c21b93c2e8ff5d2f0b01ba6f1462c3e1.de-de.xlf", "", "", "91bbbb67-038f-4402-9f8c-f2af61374c0d.c31ca651c21b93c2e8ff5d2f0b01ba6f1462c3e1.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/b3816e8e29a08168179fb144aefcdd86e306b64e/e2e/91bbbb67-038f-4402-9f8c-f2af61374c0d.md", "", "", "91bbbb67-038f-4402-9f8c-f2af61374c0d.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ce3bff8d724c7278db06d293b6a480cd211a3519/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/91bbbb67-038f-4402-9f8c-f2af61374c0d.c31ca651c21b93c2e8ff5d2f0b01ba6f1462c3e1.de-de.xlf", "", "", "91bbbb67-038f-4402-9f8c-f2af61374c0d.c31ca651c21b93c2e8ff5d2f0b01ba6f1462c3e1.de-de.xlf") | Out-Null

$wsDe.Range("E2").Value = $deHandoffDt
$wsDe.Range("H2").Value = $deHandbackDt
